$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.026182
$ws.Range("H2").Value = 0.078546
$ws.Range("I2").Value = 0.02060098669457318
$ws.Range("J2").Value = 0.02060098669457318
$ws.Range("M2").Value = 3.281109666666667
$ws.Range("N2").Value = 9.843329000000001
$ws.Range("O2").Value = 0.2779739143628921
$ws.Range("P2").Value = 0.2779739143628921
$ws.Range("Q2").Value = 0.08590601329266667
$ws.Range("R2").Value = 0.773154119634
$ws.Range("S2").Value = 0.005726536911228365
$ws.Range("T2").Value = 0.005726536911228364
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.026182
$ws.Range("H3").Value = 0.078546
$ws.Range("I3").Value = 0.02060098669457318
$ws.Range("J3").Value = 0.02060098669457318
$ws.Range("M3").Value = 6.153936333333334
$ws.Range("O3").Value = 0.5213583040808726
$ws.Range("P3").Value = 0.5213583040808725
$ws.Range("Q3").Value = 0.1611223610793334
$ws.Range("R3").Value = 1.450101249714
$ws.Range("S3").Value = 0.01074049548547529
$ws.Range("T3").Value = 0.01074049548547529
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.026182
$ws.Range("H4").Value = 0.078546
$ws.Range("I4").Value = 0.02060098669457318
$ws.Range("J4").Value = 0.02060098669457318
$ws.Range("O4").Value = 0.2006677815562353
$ws.Range("P4").Value = 0.2006677815562353
$ws.Range("Q4").Value = 0.06201506047533333
$ws.Range("R4").Value = 0.5581355442780001
$ws.Range("S4").Value = 0.004133954297869522
$ws.Range("T4").Value = 0.00413395429786952
$ws.Range("I5").Value = 0.9231010325934437
$ws.Range("J5").Value = 0.9231010325934434
$ws.Range("M5").Value = 3.281109666666667
$ws.Range("N5").Value = 9.843329000000001
$ws.Range("O5").Value = 0.2779739143628921
$ws.Range("P5").Value = 0.2779739143628921
$ws.Range("Q5").Value = 3.84932677022389
$ws.Range("R5").Value = 34.64394093201501
$ws.Range("S5").Value = 0.2565980073824272
$ws.Range("T5").Value = 0.2565980073824271
$ws.Range("I6").Value = 0.9231010325934437
$ws.Range("J6").Value = 0.9231010325934434
$ws.Range("M6").Value = 6.153936333333334
$ws.Range("O6").Value = 0.5213583040808726
$ws.Range("P6").Value = 0.5213583040808725
$ws.Range("Q6").Value = 7.219664770979446
$ws.Range("R6").Value = 64.97698293881501
$ws.Range("S6").Value = 0.4812663888482201
$ws.Range("T6").Value = 0.4812663888482199
$ws.Range("I7").Value = 0.9231010325934437
$ws.Range("J7").Value = 0.9231010325934434
$ws.Range("O7").Value = 0.2006677815562353
$ws.Range("P7").Value = 0.2006677815562353
$ws.Range("S7").Value = 0.1852366363627964
$ws.Range("T7").Value = 0.1852366363627964
$ws.Range("I8").Value = 0.05629798071198328
$ws.Range("J8").Value = 0.05629798071198327
$ws.Range("M8").Value = 3.281109666666667
$ws.Range("N8").Value = 9.843329000000001
$ws.Range("O8").Value = 0.2779739143628921
$ws.Range("P8").Value = 0.2779739143628921
$ws.Range("Q8").Value = 0.2347623029467778
$ws.Range("R8").Value = 2.112860726521
$ws.Range("S8").Value = 0.01564937006923659
$ws.Range("T8").Value = 0.01564937006923659
$ws.Range("I9").Value = 0.05629798071198328
$ws.Range("J9").Value = 0.05629798071198327
$ws.Range("M9").Value = 6.153936333333334
$ws.Range("O9").Value = 0.5213583040808726
$ws.Range("P9").Value = 0.5213583040808725
$ws.Range("Q9").Value = 0.440312093337889
$ws.Range("R9").Value = 3.962808840041001
$ws.Range("S9").Value = 0.02935141974717728
$ws.Range("T9").Value = 0.02935141974717726
$ws.Range("I10").Value = 0.05629798071198328
$ws.Range("J10").Value = 0.05629798071198327
$ws.Range("O10").Value = 0.2006677815562353
$ws.Range("P10").Value = 0.2006677815562353
$ws.Range("S10").Value = 0.01129719089556941
$ws.Range("T10").Value = 0.01129719089556941
